# Scheduled runner update: refresh computed market/profit columns (H:N)
# across the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Sheets.Item("ALC")

$ws.Range("H62").Value = 3949.1
$ws.Range("I62").Value = 2941.4285
$ws.Range("J62").Value = 6300.3335
$ws.Range("K62").Value = 2941.4285
$ws.Range("L62").Value = 6300.3335
$ws.Range("M62").Value = -2317.4285
$ws.Range("N62").Value = -7548.3335

$ws.Range("H65").Value = 3949.1
$ws.Range("I65").Value = 2941.4285
$ws.Range("J65").Value = 6300.3335
$ws.Range("K65").Value = 14707.1425
$ws.Range("L65").Value = 31501.6675
$ws.Range("M65").Value = -11587.1425
$ws.Range("N65").Value = -37741.6675

$ws.Range("H132").Value = 2532.182
$ws.Range("I132").Value = 2532.182
$ws.Range("K132").Value = 7596.545999999999
$ws.Range("M132").Value = -5066.545999999999

$ws.Range("H141").Value = 3092.1875
$ws.Range("I141").Value = 2792.3572
$ws.Range("J141").Value = 5191
$ws.Range("K141").Value = 8377.071599999999
$ws.Range("L141").Value = 15573
$ws.Range("M141").Value = -3197.071599999999
$ws.Range("N141").Value = -25933

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Sheets.Item("ARM")

$ws.Range("H5").Value = 109.42857
$ws.Range("I5").Value = 131.4
$ws.Range("J5").Value = 54.5
$ws.Range("K5").Value = 131.4
$ws.Range("L5").Value = 54.5
$ws.Range("M5").Value = -19.40000000000001
$ws.Range("N5").Value = -278.5

$ws.Range("H122").Value = 2138.7144
$ws.Range("I122").Value = 2138.7144
$ws.Range("K122").Value = 6416.1432
$ws.Range("M122").Value = -3966.1432

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Sheets.Item("BSM")

$ws.Range("H4").Value = 109.42857
$ws.Range("I4").Value = 131.4
$ws.Range("J4").Value = 54.5
$ws.Range("K4").Value = 131.4
$ws.Range("L4").Value = 54.5
$ws.Range("M4").Value = -16.40000000000001
$ws.Range("N4").Value = -284.5

$ws.Range("H11").Value = 327.33334
$ws.Range("I11").Value = 241
$ws.Range("K11").Value = 241
$ws.Range("M11").Value = -101

$ws.Range("H20").Value = 4867.4
$ws.Range("I20").Value = 3584.25
$ws.Range("K20").Value = 3584.25
$ws.Range("M20").Value = -3337.25

$ws.Range("H94").Value = 1869.4445
$ws.Range("I94").Value = 2102.6667
$ws.Range("K94").Value = 2102.6667
$ws.Range("M94").Value = -1651.6667

$ws.Range("H107").Value = 811
$ws.Range("I107").Value = 811
$ws.Range("K107").Value = 811
$ws.Range("M107").Value = 1109

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H134").Value = 8684.27
$ws.Range("I134").Value = 10825.412
$ws.Range("J134").Value = 4639.8887
$ws.Range("K134").Value = 32476.236
$ws.Range("L134").Value = 13919.6661
$ws.Range("M134").Value = -29941.236
$ws.Range("N134").Value = -18989.6661

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Sheets.Item("CRP")

$ws.Range("H7").Value = 148.28572
$ws.Range("I7").Value = 76.59999999999999
$ws.Range("J7").Value = 327.5
$ws.Range("K7").Value = 76.59999999999999
$ws.Range("L7").Value = 327.5
$ws.Range("M7").Value = 36.40000000000001
$ws.Range("N7").Value = -553.5

$ws.Range("H62").Value = 3001.3333
$ws.Range("I62").Value = 1002.5
$ws.Range("J62").Value = 6999
$ws.Range("K62").Value = 1002.5
$ws.Range("L62").Value = 6999
$ws.Range("M62").Value = -378.5
$ws.Range("N62").Value = -8247

$ws.Range("H65").Value = 3001.3333
$ws.Range("I65").Value = 1002.5
$ws.Range("J65").Value = 6999
$ws.Range("K65").Value = 5012.5
$ws.Range("L65").Value = 34995
$ws.Range("M65").Value = -1892.5
$ws.Range("N65").Value = -41235

$ws.Range("H74").Value = 39728
$ws.Range("J74").Value = 39997.5
$ws.Range("L74").Value = 39997.5
$ws.Range("N74").Value = -41745.5

$ws.Range("H77").Value = 39728
$ws.Range("J77").Value = 39997.5
$ws.Range("L77").Value = 119992.5
$ws.Range("N77").Value = -128728.5

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 2239
$ws.Range("I132").Value = 1987
$ws.Range("K132").Value = 5961
$ws.Range("M132").Value = -3431

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Sheets.Item("GSM")

$ws.Range("H70").Value = 6183.6665
$ws.Range("I70").Value = 4366.3335
$ws.Range("K70").Value = 4366.3335
$ws.Range("M70").Value = -4096.3335

$ws.Range("H73").Value = 6183.6665
$ws.Range("I73").Value = 4366.3335
$ws.Range("K73").Value = 4366.3335
$ws.Range("M73").Value = -3430.3335

$ws.Range("H80").Value = 2725.182
$ws.Range("I80").Value = 1697
$ws.Range("J80").Value = 2953.6667
$ws.Range("K80").Value = 1697
$ws.Range("L80").Value = 2953.6667
$ws.Range("M80").Value = -699
$ws.Range("N80").Value = -4949.6667

$ws.Range("H83").Value = 2725.182
$ws.Range("I83").Value = 1697
$ws.Range("J83").Value = 2953.6667
$ws.Range("K83").Value = 8485
$ws.Range("L83").Value = 14768.3335
$ws.Range("M83").Value = -3493
$ws.Range("N83").Value = -24752.3335

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Sheets.Item("LTW")

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws.Range("H62").Value = 49999.668
$ws.Range("J62").Value = 49999.668
$ws.Range("L62").Value = 49999.668
$ws.Range("N62").Value = -51247.668

$ws.Range("H65").Value = 49999.668
$ws.Range("J65").Value = 49999.668
$ws.Range("L65").Value = 149999.004
$ws.Range("N65").Value = -156239.004

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Sheets.Item("WVR")

$ws.Range("H68").Value = 14750
$ws.Range("J68").Value = 14750
$ws.Range("L68").Value = 14750
$ws.Range("N68").Value = -16372

$ws.Range("H71").Value = 14750
$ws.Range("J71").Value = 14750
$ws.Range("L71").Value = 44250
$ws.Range("N71").Value = -52362

$ws.Range("H136").Value = 1947.5
$ws.Range("I136").Value = 2014.7778
$ws.Range("J136").Value = 1745.6666
$ws.Range("K136").Value = 6044.3334
$ws.Range("L136").Value = 5236.9998
$ws.Range("M136").Value = -3494.3334
$ws.Range("N136").Value = -10336.9998
